$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10418609
$ws.Range("J17").Value = 10418609
$ws.Range("L17").Value = 31255827
$ws.Range("N17").Value = -31256163

$ws.Range("H70").Value = 1232.4445
$ws.Range("J70").Value = 998.25
$ws.Range("L70").Value = 2994.75
$ws.Range("N70").Value = -3534.75

$ws.Range("H73").Value = 1232.4445
$ws.Range("J73").Value = 998.25
$ws.Range("L73").Value = 2994.75
$ws.Range("N73").Value = -4866.75

$ws.Range("H80").Value = 522.8461
$ws.Range("I80").Value = 533.25
$ws.Range("K80").Value = 1599.75
$ws.Range("M80").Value = -601.75

$ws.Range("H83").Value = 522.8461
$ws.Range("I83").Value = 533.25
$ws.Range("K83").Value = 4799.25
$ws.Range("M83").Value = 192.75

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 9092055
$ws.Range("J88").Value = 1311.2222
$ws.Range("L88").Value = 1311.2222
$ws.Range("N88").Value = -2123.2222

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 9092055
$ws.Range("J91").Value = 1311.2222
$ws.Range("L91").Value = 1311.2222
$ws.Range("N91").Value = -4119.2222

$ws.Range("H137").Value = 4555873
$ws.Range("I137").Value = 6250950
$ws.Range("K137").Value = 18752850
$ws.Range("M137").Value = -18750300

$ws.Range("H138").Value = 2909.904
$ws.Range("I138").Value = 2970.3
$ws.Range("K138").Value = 8910.900000000001
$ws.Range("M138").Value = -3770.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 18500
$ws.Range("J27").Value = 18500
$ws.Range("L27").Value = 18500
$ws.Range("N27").Value = -18868

$ws.Range("H61").Value = 3216.0857
$ws.Range("I61").Value = 2420.926
$ws.Range("J61").Value = 5899.75
$ws.Range("K61").Value = 2420.926
$ws.Range("L61").Value = 5899.75
$ws.Range("M61").Value = -2208.926
$ws.Range("N61").Value = -6323.75

$ws.Range("H74").Value = 267551.25
$ws.Range("I74").Value = 619195.9
$ws.Range("J74").Value = 3817.75
$ws.Range("K74").Value = 619195.9
$ws.Range("L74").Value = 3817.75
$ws.Range("M74").Value = -618321.9
$ws.Range("N74").Value = -5565.75

$ws.Range("H77").Value = 267551.25
$ws.Range("I77").Value = 619195.9
$ws.Range("J77").Value = 3817.75
$ws.Range("K77").Value = 3095979.5
$ws.Range("L77").Value = 19088.75
$ws.Range("M77").Value = -3091611.5
$ws.Range("N77").Value = -27824.75

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws.Range("H136").Value = 3216.0857
$ws.Range("I136").Value = 2420.926
$ws.Range("J136").Value = 5899.75
$ws.Range("K136").Value = 7262.778
$ws.Range("L136").Value = 17699.25
$ws.Range("M136").Value = -4712.778
$ws.Range("N136").Value = -22799.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4041.5
$ws.Range("I134").Value = 3913.7273
$ws.Range("J134").Value = 4582.077
$ws.Range("K134").Value = 11741.1819
$ws.Range("L134").Value = 13746.231
$ws.Range("M134").Value = -9206.1819
$ws.Range("N134").Value = -18816.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3876.96
$ws.Range("I31").Value = 2866.8572
$ws.Range("K31").Value = 2866.8572
$ws.Range("M31").Value = -2571.8572

$ws.Range("H34").Value = 3876.96
$ws.Range("I34").Value = 2866.8572
$ws.Range("K34").Value = 2866.8572
$ws.Range("M34").Value = -2664.8572

$ws.Range("H58").Value = 3043.9644
$ws.Range("I58").Value = 2497.5
$ws.Range("K58").Value = 2497.5
$ws.Range("M58").Value = -2294.5

$ws.Range("H62").Value = 16682815
$ws.Range("I62").Value = 33341966
$ws.Range("J62").Value = 23662.666
$ws.Range("K62").Value = 33341966
$ws.Range("L62").Value = 23662.666
$ws.Range("M62").Value = -33341342
$ws.Range("N62").Value = -24910.666

$ws.Range("H65").Value = 16682815
$ws.Range("I65").Value = 33341966
$ws.Range("J65").Value = 23662.666
$ws.Range("K65").Value = 166709830
$ws.Range("L65").Value = 118313.33
$ws.Range("M65").Value = -166706710
$ws.Range("N65").Value = -124553.33

$ws.Range("H94").Value = 2196.4375
$ws.Range("J94").Value = 2332.3333
$ws.Range("L94").Value = 2332.3333
$ws.Range("N94").Value = -3234.3333

$ws.Range("H132").Value = 3839.682
$ws.Range("I132").Value = 4157.625
$ws.Range("J132").Value = 3658
$ws.Range("K132").Value = 12472.875
$ws.Range("L132").Value = 10974
$ws.Range("M132").Value = -9942.875
$ws.Range("N132").Value = -16034

$ws.Range("H136").Value = 3043.9644
$ws.Range("I136").Value = 2497.5
$ws.Range("K136").Value = 7492.5
$ws.Range("M136").Value = -4942.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 503000
$ws.Range("I70").Value = 503000
$ws.Range("K70").Value = 503000
$ws.Range("M70").Value = -502730

$ws.Range("H73").Value = 503000
$ws.Range("I73").Value = 503000
$ws.Range("K73").Value = 503000
$ws.Range("M73").Value = -502064

$ws.Range("H97").Value = 4034.4
$ws.Range("I97").Value = 4511.25
$ws.Range("J97").Value = 3716.5
$ws.Range("K97").Value = 4511.25
$ws.Range("L97").Value = 3716.5
$ws.Range("M97").Value = -4015.25
$ws.Range("N97").Value = -4708.5

$ws.Range("H102").Value = 1859.2667
$ws.Range("I102").Value = 1221.1111
$ws.Range("J102").Value = 2816.5
$ws.Range("K102").Value = 1221.1111
$ws.Range("L102").Value = 2816.5
$ws.Range("M102").Value = 400.8888999999999
$ws.Range("N102").Value = -6060.5

$ws.Range("H126").Value = 19850.334
$ws.Range("I126").Value = 17555
$ws.Range("K126").Value = 52665
$ws.Range("M126").Value = -50195

$ws.Range("H132").Value = 2787.375
$ws.Range("I132").Value = 1383.1666
$ws.Range("K132").Value = 4149.4998
$ws.Range("M132").Value = -1619.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 31379.621
$ws.Range("I40").Value = 38006.066
$ws.Range("K40").Value = 38006.066
$ws.Range("M40").Value = -37870.066

$ws.Range("H55").Value = 334.42856
$ws.Range("I55").Value = 200.25
$ws.Range("K55").Value = 200.25
$ws.Range("M55").Value = -27.25

$ws.Range("H61").Value = 2611.158
$ws.Range("I61").Value = 2518.0588
$ws.Range("J61").Value = 3402.5
$ws.Range("K61").Value = 2518.0588
$ws.Range("L61").Value = 3402.5
$ws.Range("M61").Value = -2316.0588
$ws.Range("N61").Value = -3806.5

$ws.Range("H82").Value = 593.2857
$ws.Range("I82").Value = 584
$ws.Range("J82").Value = 605.6667
$ws.Range("K82").Value = 584
$ws.Range("L82").Value = 605.6667
$ws.Range("M82").Value = -223
$ws.Range("N82").Value = -1327.6667

$ws.Range("H85").Value = 593.2857
$ws.Range("I85").Value = 584
$ws.Range("J85").Value = 605.6667
$ws.Range("K85").Value = 584
$ws.Range("L85").Value = 605.6667
$ws.Range("M85").Value = 664
$ws.Range("N85").Value = -3101.6667

$ws.Range("H113").Value = 2611.158
$ws.Range("I113").Value = 2518.0588
$ws.Range("J113").Value = 3402.5
$ws.Range("K113").Value = 2518.0588
$ws.Range("L113").Value = 3402.5
$ws.Range("M113").Value = -348.0587999999998
$ws.Range("N113").Value = -7742.5

$ws.Range("H122").Value = 4020.5
$ws.Range("I122").Value = 4200.778
$ws.Range("K122").Value = 12602.334
$ws.Range("M122").Value = -10152.334

$ws.Range("H132").Value = 5223.4585
$ws.Range("I132").Value = 3025.8572
$ws.Range("J132").Value = 8300.1
$ws.Range("K132").Value = 9077.571599999999
$ws.Range("L132").Value = 24900.3
$ws.Range("M132").Value = -6547.571599999999
$ws.Range("N132").Value = -29960.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H96").Value = 4182.7896
$ws.Range("I96").Value = 4553.5454
$ws.Range("J96").Value = 3673
$ws.Range("K96").Value = 4553.5454
$ws.Range("L96").Value = 3673
$ws.Range("M96").Value = -3180.5454
$ws.Range("N96").Value = -6419

$ws.Range("H132").Value = 8250
$ws.Range("I132").Value = 11250
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 33750
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -31220
$ws.Range("N132").Value = -20810
